$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I2: replace the =TRUE() boolean formula with the literal text "TRUE"
$c2 = $ws.Cells.Item(2, 9)
$c2.NumberFormat = "@"
$c2.Formula = "=""TRUE"""
$c2.Copy()
$c2.PasteSpecial(-4163)

# I3: replace the =TRUE() boolean formula with the literal text "TRUE"
$c3 = $ws.Cells.Item(3, 9)
$c3.NumberFormat = "@"
$c3.Formula = "=""TRUE"""
$c3.Copy()
$c3.PasteSpecial(-4163)

# Move / collapse the selection onto I3 (matches the post-edit selection state)
$ws.Range("I3").Select()
